# Update cryptos list: price (D) and volume/1h (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.862.10"
$ws.Range("E2").Value = "  -4.05%  "

$ws.Range("D3").Value = "3.467.27"
$ws.Range("E3").Value = "  -3.96%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.83%  "

$ws.Range("D7").Value = "3.466.13"
$ws.Range("E7").Value = "  -3.97%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("E10").Value = "  -4.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("E12").Value = "  -3.99%  "

$ws.Range("E13").Value = "  -4.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.93%  "

$ws.Range("D15").Value = "4.052.26"
$ws.Range("E15").Value = "  -3.97%  "

$ws.Range("D16").Value = "3.467.47"
$ws.Range("E16").Value = "  -3.96%  "

$ws.Range("D17").Value = "66.928.82"
$ws.Range("E17").Value = "  -3.73%  "

$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.41%  "

$ws.Range("E20").Value = "  -4.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "440.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.25%  "

$ws.Range("E23").Value = "  -5.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Value = "3.601.96"
$ws.Range("E26").Value = "  -4.07%  "

$ws.Range("E27").Value = "  -10.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.15%  "

$ws.Range("E32").Value = "  -3.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.63%  "

$ws.Range("D36").Value = "3.457.86"
$ws.Range("E36").Value = "  -4.13%  "

$ws.Range("E37").Value = "  -8.04%  "

$ws.Range("E38").Value = "  -6.65%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.87%  "

$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("E43").Value = "  -11.40%  "

$ws.Range("E44").Value = "  -5.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.886"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("E48").Value = "  -10.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.988"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.71%  "
